# Account Application.docx backend-fields update
#
# 1) The list item reading "Expense Category" (the one immediately
#    preceding the "Category Name" item) becomes just "Category".
# 2) A new list item "Category Type" is inserted right after
#    "Category Name" (before "Category Description"), using the same
#    list formatting (ListParagraph style / numId 6) as its neighbours.

$d = $word.ActiveDocument
$paras = $d.Paragraphs
$count = $paras.Count

function Get-ParaText($para) {
    return $para.Range.Text.Trim()
}

# --- Step 1: rename the "Expense Category" list item to "Category" ---
$expenseCategoryIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $text = Get-ParaText $paras.Item($i)
    if ($text -eq "Expense Category") {
        if ($i -lt $count -and (Get-ParaText $paras.Item($i + 1)) -eq "Category Name") {
            $expenseCategoryIndex = $i
            break
        }
    }
}

if ($expenseCategoryIndex -eq -1) {
    throw "Could not locate the 'Expense Category' list item preceding 'Category Name'."
}

$paras.Item($expenseCategoryIndex).Range.Text = "Category"

# --- Step 2: insert "Category Type" right after "Category Name" ---
# Recompute the paragraphs collection / count since the document changed.
$paras = $d.Paragraphs
$count = $paras.Count

$categoryNameIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $text = Get-ParaText $paras.Item($i)
    if ($text -eq "Category Name") {
        $categoryNameIndex = $i
        break
    }
}

if ($categoryNameIndex -eq -1) {
    throw "Could not locate the 'Category Name' list item."
}

$categoryNamePara = $paras.Item($categoryNameIndex)
$categoryNamePara.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$newPara = $paras.Item($categoryNameIndex + 1)
$newPara.Range.Text = "Category Type"
